$d = $word.ActiveDocument

# 1. "...appeared in Court for sentencing on June 18, 2022." -> June 19, 2022
#    (entirely inside one run, so this cannot disturb neighbouring runs)
$d.Content.Find.Execute("appeared in Court for sentencing on June 18, 2022.", $true, $false, $false, $false, $false, $true, 1, $false, "appeared in Court for sentencing on June 19, 2022.", 2)

# 2. "...shall pay the fines and costs in full by " + bold "June 18, 2022" + "."
#    The date here is its own (bold) run sandwiched between "in full by " and
#    ".", so anchor on the unique preceding text and only touch the date
#    itself -- this keeps the bold run intact instead of swallowing its
#    differently-formatted neighbours.
$anchor2 = $d.Content
$anchor2.Find.Execute("in full by ") | Out-Null
$anchor2.Collapse(0)
$dateRun2 = $d.Range($anchor2.Start, $anchor2.Start + 13)   # "June 18, 2022" has 13 chars
$dateRun2.Text = "June 19, 2022"

# 3. "...show proof of completion ... on or before August 17, 2022." -> August 18, 2022
$d.Content.Find.Execute("Office of Community Control on or before August 17, 2022.", $true, $false, $false, $false, $false, $true, 1, $false, "Office of Community Control on or before August 18, 2022.", 2)

# 4. "Defendant's driving license is suspended from June 18, 2022" -> June 19, 2022
$d.Content.Find.Execute([char]0x2019 + "s driving license is suspended from June 18, 2022", $true, $false, $false, $false, $false, $true, 1, $false, [char]0x2019 + "s driving license is suspended from June 19, 2022", 2)
